# Applies the "updated Gender xpath to Sex xpath for canine" edit:
# the test-automation re-run now also emits a per-query "Message" log
# sheet plus a (currently empty) "StatOutput" sheet and its own
# "Message" log sheet, in addition to the existing CypherOutput/Message
# sheets.

$wb = $excel.ActiveWorkbook

# Common message-log text reused across the new "Message" sheets -
# identical to what is already on the existing "Message" sheet, except
# the output path now points at the relocated "Script" subfolder.
$neo4jUrlLabel  = "Neo4j_URL:"
$neo4jUrl       = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userLabel      = "User_name:"
$userName       = "neo4j"
$pwdLabel       = "PWD:"
$pwdValue       = "icdcDBneo4j0"
$cypherLabel    = "Cypher:"
$cypherQuery    = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.sex IN ['Castrated male'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"
$outputLabel    = "Output:"
$outputPath     = "C:\Users\radhakrishnang2\Desktop\Script\Commons_Automation\OutputFiles\TC01_Canine_Filter_Gender-CastratedMale_Neo4jData.xlsx"
$emptyStatus    = ""
$cypherEmptyMsg = "Cypher query should not be an empty string"

# ---------------------------------------------------------------------
# Sheet 3: CypherOutput_Message
# (Worksheets.Add(Before, After) - pass After the current last sheet so
# the new sheet lands at the end, keeping tab order stable without ever
# relying on .Move(), which leaves stale sheet references behind.)
# ---------------------------------------------------------------------
$cypherMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$cypherMsg.Name = "CypherOutput_Message"

$cypherMsg.Range("A1").Value  = $neo4jUrlLabel
$cypherMsg.Range("A2").Value  = $neo4jUrl
$cypherMsg.Range("A3").Value  = $userLabel
$cypherMsg.Range("A4").Value  = $userName
$cypherMsg.Range("A5").Value  = $pwdLabel
$cypherMsg.Range("A6").Value  = $pwdValue
$cypherMsg.Range("A7").Value  = $cypherLabel
$cypherMsg.Range("A8").Value  = $cypherQuery
$cypherMsg.Range("A9").Value  = $outputLabel
$cypherMsg.Range("A10").Value = $outputPath

# ---------------------------------------------------------------------
# Sheet 4: StatOutput (empty result sheet)
# ---------------------------------------------------------------------
$statOutput = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$statOutput.Name = "StatOutput"

# ---------------------------------------------------------------------
# Sheet 5: StatOutput_Message
# ---------------------------------------------------------------------
$statMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$statMsg.Name = "StatOutput_Message"

$statMsg.Range("A1").Value  = $neo4jUrlLabel
$statMsg.Range("A2").Value  = $neo4jUrl
$statMsg.Range("A3").Value  = $userLabel
$statMsg.Range("A4").Value  = $userName
$statMsg.Range("A5").Value  = $pwdLabel
$statMsg.Range("A6").Value  = $pwdValue
$statMsg.Range("A7").Value  = $cypherLabel
$statMsg.Range("A8").Value  = $cypherQuery
$statMsg.Range("A9").Value  = $outputLabel
$statMsg.Range("A10").Value = $outputPath
$statMsg.Range("A11").Value = $cypherEmptyMsg
$statMsg.Range("A12").Value = $neo4jUrlLabel
$statMsg.Range("A13").Value = $neo4jUrl
$statMsg.Range("A14").Value = $userLabel
$statMsg.Range("A15").Value = $userName
$statMsg.Range("A16").Value = $pwdLabel
$statMsg.Range("A17").Value = $pwdValue
$statMsg.Range("A18").Value = $cypherLabel
$statMsg.Range("A19").Value = $emptyStatus
$statMsg.Range("A20").Value = $outputLabel
$statMsg.Range("A21").Value = $outputPath

# Restore the originally-active sheet/tab (CypherOutput) so activeTab
# stays 0 instead of pointing at the newly-inserted last sheet.
$wb.Worksheets.Item(1).Activate()
